$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Heading, level 1*") {
        $p.Range.Delete()
        break
    }
}
